# Update the Files-tab SQL query in cell B5: drop the file_source column
# (and its File Source output column) from the SELECT list, matching the
# updated query used for the phs002371 test-case regeneration.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B5").Value = 'WITH file_data AS (
    SELECT 
        file_name, 
        data_category,
        file_type, 
        file_size,
        file_access,  
        file_description,
        "sample.id",
        ''Sequencing'' AS file_source
    FROM df_sequencing_file

    UNION

    SELECT 
        file_name, 
        data_category,
        file_type, 
        file_size,
        file_access,  
        file_description,
        "sample.id",
        ''Pathology'' AS file_source
    FROM df_pathology_file
)

SELECT DISTINCT
    fd.file_name AS "File Name",
    fd.data_category AS "Data Category",
    COALESCE(fd.file_description, '''') AS "File Description",
    fd.file_type AS "File Type",
    CASE
        WHEN fd.file_size >= 1024 * 1024 * 1024 THEN
            CASE 
                WHEN ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 2) = CAST(ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT)
                THEN CAST(CAST(ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) AS TEXT) || '' GB''
                ELSE ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 2) || '' GB''
            END
        WHEN fd.file_size >= 1024 * 1024 THEN
            CASE 
                WHEN ROUND(fd.file_size / (1024.0 * 1024.0), 2) = CAST(ROUND(fd.file_size / (1024.0 * 1024.0), 0) AS INT)
                THEN CAST(CAST(ROUND(fd.file_size / (1024.0 * 1024.0), 0) AS INT) AS TEXT) || '' MB''
                ELSE ROUND(fd.file_size / (1024.0 * 1024.0), 2) || '' MB''
            END
        WHEN fd.file_size >= 1024 THEN
            CASE 
                WHEN ROUND(fd.file_size / 1024.0, 2) = CAST(ROUND(fd.file_size / 1024.0, 0) AS INT)
                THEN CAST(CAST(ROUND(fd.file_size / 1024.0, 0) AS INT) AS TEXT) || '' KB''
                ELSE ROUND(fd.file_size / 1024.0, 2) || '' KB''
            END
        ELSE 
            CASE 
                WHEN ROUND(fd.file_size, 2) = CAST(ROUND(fd.file_size, 0) AS INT)
                THEN CAST(CAST(ROUND(fd.file_size, 0) AS INT) AS TEXT) || '' Bytes''
                ELSE ROUND(fd.file_size, 2) || '' Bytes''
            END
    END AS "File Size",
    fd.file_access AS "File Access",
    std.dbgap_accession AS "Study ID",
    prt.participant_id AS "Participant ID",
    smp.sample_id AS "Sample ID"
  FROM 
    df_study std
LEFT JOIN df_participant prt ON std.id = prt."study.id"
LEFT JOIN df_sample smp ON prt.id = smp."participant.id"
JOIN file_data fd ON smp.id = fd."sample.id"
LEFT JOIN df_diagnosis dgn ON prt.id = dgn."participant.id"
WHERE 
    std.dbgap_accession = ''phs000720''
    AND prt.sex_at_birth = ''Female''
    AND prt.race LIKE ''%White%''
    AND dgn.disease_phase = ''Relapse''
ORDER BY fd.file_name 
LIMIT 100;'
